# untypedNamedBlanks.xlsx - fix failing atomic tabOTTR tests due to breaking change.
#
# The three data cells A7:A9 on Sheet1 get "rotated": the literal number 1
# moves up from A8 to A7, and the two text labels ("auto"/"data") shift down
# into A8/A9 respectively. A7/A9 also drop back to the sheet's base
# (un-applied-font) style, while A8 keeps the bold/explicit-font style that
# used to live on A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 1
$ws.Range("A7").NumberFormat = "General"

$ws.Range("A8").Value = "auto"

$ws.Range("A9").Value = "data"
$ws.Range("A9").NumberFormat = "General"

# Move the active selection from A10 to A8, matching the saved view state.
$ws.Range("A8").Select()

# Best-effort: nudge the window's tab-ratio setting to match the saved view.
$excel.ActiveWindow.TabRatio = 990
